$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the "trump" player's data (merged into row 2, replacing loclexyz99)
$ws.Range("A2").Value = "trump"
$ws.Range("B2").Value = "f"
$ws.Range("C2").Value = "trump.png"
$ws.Range("D2").Value = "a"
$ws.Range("E2").Value = 45
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 0.5

# Delete row 3 entirely (the old "trump" row)
$ws.Rows.Item(3).Delete()

# Update selection to D4 as per the diff
$ws.Range("D4").Select()
